# Update countries & provincias Spain
#
# This mirrors the 15 Abr 2020 17:52 -> 18:22 data refresh:
#  - A handful of countries (India, Irak, Guinea Ecuatorial) moved up two
#    rows in the sorted ("Casos totales" desc) list because their totals
#    grew past the countries that used to sit just above them; those
#    displaced countries (Israel/Suecia, Kuwait/Estonia, Bahamas/Guyana/
#    Zambia/Macao/Guinea-Bisau) shift down a row to make room.
#  - The statistic columns (B:H) are refreshed for every row whose numbers
#    actually changed between the two snapshots.
#  - The "Datos actualizados ..." banner in A1 gets the new timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 18:22"

function Set-RowValues {
    param($Sheet, [int]$Row, [string]$Pais, $Values)

    $Sheet.Cells.Item($Row, 1).Value = $Pais

    $n = $Values.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCol = 2   # column B
    $endCol = $startCol + $n - 1
    $rng = $Sheet.Range($Sheet.Cells.Item($Row, $startCol), $Sheet.Cells.Item($Row, $endCol))
    $rng.Value = $arr
}

# Row => final (Pais, B..H) state after the refresh/re-sort.
# Rows not listed here are unchanged by this edit.
Set-RowValues $ws 4   "Estados Unidos"    @(618856, 4970, 40271, 552252, 13473, 286, 26333)
Set-RowValues $ws 6   "Italia"            @(165155, 2667, 38092, 105418, 3079, 578, 21645)
Set-RowValues $ws 8   "Alemania"          @(132747, 537, 72600, 56555, 4288, 97, 3592)
Set-RowValues $ws 17  "Brasil"            @(26113, 851, 14026, 10497, 296, 58, 1590)
Set-RowValues $ws 21  "India"             @(12320, 833, 1432, 10483, 0, 12, 405)
Set-RowValues $ws 22  "Israel"            @(12200, 154, 2309, 9765, 176, 3, 126)
Set-RowValues $ws 23  "Suecia"            @(11927, 482, 381, 10343, 954, 170, 1203)
Set-RowValues $ws 58  "Argelia"           @(2160, 90, 708, 1116, 60, 10, 336)
Set-RowValues $ws 65  "Irak"              @(1415, 15, 812, 524, 0, 1, 79)
Set-RowValues $ws 66  "Kuwait"            @(1405, 50, 206, 1196, 31, 0, 3)
Set-RowValues $ws 67  "Estonia"           @(1400, 27, 117, 1248, 10, 4, 35)
Set-RowValues $ws 112 "Georgia"           @(306, 6, 71, 232, 6, 0, 3)
Set-RowValues $ws 117 "Sri Lanka"         @(237, 4, 63, 167, 1, 0, 7)
Set-RowValues $ws 154 "Guinea Ecuatorial" @(51, 10, 4, 47, 0, 0, 0)
Set-RowValues $ws 155 "Bahamas"           @(49, 0, 6, 35, 1, 0, 8)
Set-RowValues $ws 156 "Guyana"            @(48, 1, 8, 34, 5, 0, 6)
Set-RowValues $ws 157 "Zambia"            @(48, 3, 30, 16, 1, 0, 2)
Set-RowValues $ws 158 "Macao"             @(45, 0, 15, 30, 1, 0, 0)
Set-RowValues $ws 159 "Guinea-Bisau"      @(43, 0, 0, 43, 0, 0, 0)

$wb.Save()
